# minor fix on date
# Cell A3 ("Project review meeting with TA" row) had a stray trailing slash
# in its date string: " 11/16-11/19/" -> " 11/16-11/19"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = " 11/16-11/19"
